$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 37037400
$ws.Range("I28").Value = 37037400
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 37037400
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -37036915
$ws.Range("N28").Value = $null
$ws.Range("H32").Value = 100008340
$ws.Range("I32").Value = 100012500
$ws.Range("K32").Value = 100012500
$ws.Range("M32").Value = -100012174
$ws.Range("H33").Value = 762.8182
$ws.Range("I33").Value = 869.94116
$ws.Range("J33").Value = 398.6
$ws.Range("K33").Value = 869.94116
$ws.Range("L33").Value = 398.6
$ws.Range("M33").Value = -640.94116
$ws.Range("N33").Value = -856.6
$ws.Range("H53").Value = 447.9
$ws.Range("I53").Value = 563.75
$ws.Range("J53").Value = 370.66666
$ws.Range("K53").Value = 563.75
$ws.Range("L53").Value = 370.66666
$ws.Range("M53").Value = 73.25
$ws.Range("N53").Value = -1644.66666
$ws.Range("H97").Value = 1380.2858
$ws.Range("J97").Value = 1404.3
$ws.Range("L97").Value = 4212.9
$ws.Range("N97").Value = -5204.9
$ws.Range("H98").Value = 1353.7715
$ws.Range("I98").Value = 703.96155
$ws.Range("J98").Value = 3231
$ws.Range("K98").Value = 703.96155
$ws.Range("L98").Value = 3231
$ws.Range("M98").Value = 794.03845
$ws.Range("N98").Value = -6227
$ws.Range("H107").Value = 19234784
$ws.Range("I107").Value = 22730680
$ws.Range("J107").Value = 7349.5
$ws.Range("K107").Value = 22730680
$ws.Range("L107").Value = 7349.5
$ws.Range("M107").Value = -22728760
$ws.Range("N107").Value = -11189.5
$ws.Range("H112").Value = 1302.3334
$ws.Range("J112").Value = 1261.4263
$ws.Range("L112").Value = 3784.2789
$ws.Range("N112").Value = -6000.2789
$ws.Range("H116").Value = 6132.522
$ws.Range("I116").Value = 6658.077
$ws.Range("J116").Value = 5449.3
$ws.Range("K116").Value = 6658.077
$ws.Range("L116").Value = 5449.3
$ws.Range("M116").Value = -3216.077
$ws.Range("N116").Value = -12333.3
$ws.Range("H122").Value = 1353.7715
$ws.Range("I122").Value = 703.96155
$ws.Range("J122").Value = 3231
$ws.Range("K122").Value = 2111.88465
$ws.Range("L122").Value = 9693
$ws.Range("M122").Value = 338.11535
$ws.Range("N122").Value = -14593
$ws.Range("H124").Value = 99000
$ws.Range("J124").Value = 99000
$ws.Range("L124").Value = 99000
$ws.Range("N124").Value = -108820
$ws.Range("H132").Value = 5166.051
$ws.Range("I132").Value = 2426.9778
$ws.Range("K132").Value = 7280.9334
$ws.Range("M132").Value = -4750.9334
$ws.Range("H135").Value = 1178.3556
$ws.Range("J135").Value = 5122.25
$ws.Range("L135").Value = 46100.25
$ws.Range("N135").Value = -51170.25
$ws.Range("H137").Value = 4521.362
$ws.Range("I137").Value = 6304.7085
$ws.Range("J137").Value = 2660.4783
$ws.Range("K137").Value = 18914.1255
$ws.Range("L137").Value = 7981.4349
$ws.Range("M137").Value = -16364.1255
$ws.Range("N137").Value = -13081.4349
$ws.Range("H138").Value = 2234.8853
$ws.Range("I138").Value = 1407.8379
$ws.Range("J138").Value = 3509.9167
$ws.Range("K138").Value = 4223.5137
$ws.Range("L138").Value = 10529.7501
$ws.Range("M138").Value = 916.4863000000005
$ws.Range("N138").Value = -20809.7501
$ws.Range("H141").Value = 6231.7905
$ws.Range("J141").Value = 10780.167
$ws.Range("L141").Value = 32340.501
$ws.Range("N141").Value = -42700.501

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3620.37
$ws.Range("I32").Value = 3145.0205
$ws.Range("K32").Value = 3145.0205
$ws.Range("M32").Value = -2858.0205
$ws.Range("H61").Value = 2933.775
$ws.Range("I61").Value = 2744.4
$ws.Range("K61").Value = 2744.4
$ws.Range("M61").Value = -2532.4
$ws.Range("H63").Value = 3417.7446
$ws.Range("I63").Value = 2997
$ws.Range("J63").Value = 3436.4443
$ws.Range("K63").Value = 2997
$ws.Range("L63").Value = 3436.4443
$ws.Range("M63").Value = -2311
$ws.Range("N63").Value = -4808.4443
$ws.Range("H66").Value = 3417.7446
$ws.Range("I66").Value = 2997
$ws.Range("J66").Value = 3436.4443
$ws.Range("K66").Value = 14985
$ws.Range("L66").Value = 17182.2215
$ws.Range("M66").Value = -11553
$ws.Range("N66").Value = -24046.2215
$ws.Range("H74").Value = 1806.5834
$ws.Range("I74").Value = 1519.5555
$ws.Range("J74").Value = 2667.6667
$ws.Range("K74").Value = 1519.5555
$ws.Range("L74").Value = 2667.6667
$ws.Range("M74").Value = -645.5554999999999
$ws.Range("N74").Value = -4415.6667
$ws.Range("H77").Value = 1806.5834
$ws.Range("I77").Value = 1519.5555
$ws.Range("J77").Value = 2667.6667
$ws.Range("K77").Value = 7597.7775
$ws.Range("L77").Value = 13338.3335
$ws.Range("M77").Value = -3229.7775
$ws.Range("N77").Value = -22074.3335
$ws.Range("H110").Value = 2610
$ws.Range("I110").Value = 2748.75
$ws.Range("K110").Value = 2748.75
$ws.Range("M110").Value = -703.75
$ws.Range("H132").Value = 1419.25
$ws.Range("I132").Value = 1434.1111
$ws.Range("J132").Value = 1352.375
$ws.Range("K132").Value = 4302.3333
$ws.Range("L132").Value = 4057.125
$ws.Range("M132").Value = -1772.3333
$ws.Range("N132").Value = -9117.125
$ws.Range("H136").Value = 2933.775
$ws.Range("I136").Value = 2744.4
$ws.Range("K136").Value = 8233.200000000001
$ws.Range("M136").Value = -5683.200000000001

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 102631.5
$ws.Range("J35").Value = 102631.5
$ws.Range("L35").Value = 102631.5
$ws.Range("N35").Value = -103251.5
$ws.Range("H86").Value = 2071.8572
$ws.Range("I86").Value = 1946
$ws.Range("J86").Value = 2533.3333
$ws.Range("K86").Value = 1946
$ws.Range("L86").Value = 2533.3333
$ws.Range("M86").Value = -823
$ws.Range("N86").Value = -4779.3333
$ws.Range("H89").Value = 2071.8572
$ws.Range("I89").Value = 1946
$ws.Range("J89").Value = 2533.3333
$ws.Range("K89").Value = 9730
$ws.Range("L89").Value = 12666.6665
$ws.Range("M89").Value = -4114
$ws.Range("N89").Value = -23898.6665
$ws.Range("H94").Value = 689.5
$ws.Range("I94").Value = 636.03125
$ws.Range("J94").Value = 903.375
$ws.Range("K94").Value = 636.03125
$ws.Range("L94").Value = 903.375
$ws.Range("M94").Value = -185.03125
$ws.Range("N94").Value = -1805.375
$ws.Range("H99").Value = 1433.5
$ws.Range("I99").Value = 1636
$ws.Range("J99").Value = 826
$ws.Range("K99").Value = 1636
$ws.Range("L99").Value = 826
$ws.Range("M99").Value = -138
$ws.Range("N99").Value = -3822
$ws.Range("H105").Value = 2015.7778
$ws.Range("I105").Value = 1806.0714
$ws.Range("J105").Value = 2749.75
$ws.Range("K105").Value = 1806.0714
$ws.Range("L105").Value = 2749.75
$ws.Range("M105").Value = -59.07140000000004
$ws.Range("N105").Value = -6243.75
$ws.Range("H134").Value = 4874.8066
$ws.Range("J134").Value = 8260.25
$ws.Range("L134").Value = 24780.75
$ws.Range("N134").Value = -29850.75

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1208.02
$ws.Range("I58").Value = 1171.8334
$ws.Range("J58").Value = 1398
$ws.Range("K58").Value = 1171.8334
$ws.Range("L58").Value = 1398
$ws.Range("M58").Value = -968.8334
$ws.Range("N58").Value = -1804
$ws.Range("H62").Value = 142860030
$ws.Range("I62").Value = 333335550
$ws.Range("K62").Value = 333335550
$ws.Range("M62").Value = -333334926
$ws.Range("H65").Value = 142860030
$ws.Range("I65").Value = 333335550
$ws.Range("K65").Value = 1666677750
$ws.Range("M65").Value = -1666674630
$ws.Range("H86").Value = 6069079
$ws.Range("I86").Value = 13335734
$ws.Range("J86").Value = 13533.167
$ws.Range("K86").Value = 13335734
$ws.Range("L86").Value = 13533.167
$ws.Range("M86").Value = -13334611
$ws.Range("N86").Value = -15779.167
$ws.Range("H89").Value = 6069079
$ws.Range("I89").Value = 13335734
$ws.Range("J89").Value = 13533.167
$ws.Range("K89").Value = 66678670
$ws.Range("L89").Value = 67665.83499999999
$ws.Range("M89").Value = -66673054
$ws.Range("N89").Value = -78897.83499999999
$ws.Range("H132").Value = 6254.3193
$ws.Range("I132").Value = 2373.2974
$ws.Range("J132").Value = 20614.1
$ws.Range("K132").Value = 7119.8922
$ws.Range("L132").Value = 61842.3
$ws.Range("M132").Value = -4589.8922
$ws.Range("N132").Value = -66902.29999999999
$ws.Range("H134").Value = 2808.3713
$ws.Range("I134").Value = 2803.7666
$ws.Range("K134").Value = 8411.299800000001
$ws.Range("M134").Value = -5876.299800000001
$ws.Range("H136").Value = 1208.02
$ws.Range("I136").Value = 1171.8334
$ws.Range("J136").Value = 1398
$ws.Range("K136").Value = 3515.5002
$ws.Range("L136").Value = 4194
$ws.Range("M136").Value = -965.5001999999999
$ws.Range("N136").Value = -9294
$ws.Range("H141").Value = 103843
$ws.Range("I141").Value = 55222
$ws.Range("J141").Value = 136257
$ws.Range("K141").Value = 55222
$ws.Range("L141").Value = 136257
$ws.Range("M141").Value = -50042
$ws.Range("N141").Value = -146617

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 334.16666
$ws.Range("I7").Value = 326.2
$ws.Range("J7").Value = 374
$ws.Range("K7").Value = 978.5999999999999
$ws.Range("L7").Value = 1122
$ws.Range("M7").Value = -866.5999999999999
$ws.Range("N7").Value = -1346
$ws.Range("H38").Value = 257.625
$ws.Range("I38").Value = 81.57143000000001
$ws.Range("K38").Value = 244.71429
$ws.Range("M38").Value = 102.28571
$ws.Range("H50").Value = 3001.2727
$ws.Range("I50").Value = 676.6667
$ws.Range("J50").Value = 3873
$ws.Range("K50").Value = 2030.0001
$ws.Range("L50").Value = 11619
$ws.Range("M50").Value = -1549.0001
$ws.Range("N50").Value = -12581
$ws.Range("H53").Value = 3001.2727
$ws.Range("I53").Value = 676.6667
$ws.Range("J53").Value = 3873
$ws.Range("K53").Value = 2030.0001
$ws.Range("L53").Value = 11619
$ws.Range("M53").Value = -1549.0001
$ws.Range("N53").Value = -12581
$ws.Range("H59").Value = 2221.5557
$ws.Range("I59").Value = 832.5
$ws.Range("K59").Value = 2497.5
$ws.Range("M59").Value = -1957.5
$ws.Range("H69").Value = 4424.8
$ws.Range("J69").Value = 6450
$ws.Range("L69").Value = 19350
$ws.Range("N69").Value = -20972
$ws.Range("H72").Value = 4424.8
$ws.Range("J72").Value = 6450
$ws.Range("L72").Value = 58050
$ws.Range("N72").Value = -66162
$ws.Range("H74").Value = 3999
$ws.Range("I74").Value = 3999
$ws.Range("K74").Value = 11997
$ws.Range("M74").Value = -10936
$ws.Range("H77").Value = 3999
$ws.Range("I77").Value = 3999
$ws.Range("K77").Value = 35991
$ws.Range("M77").Value = -30687
$ws.Range("H92").Value = 610.7619
$ws.Range("I92").Value = 569.0909
$ws.Range("J92").Value = 656.6
$ws.Range("K92").Value = 1707.2727
$ws.Range("L92").Value = 1969.8
$ws.Range("M92").Value = -459.2727
$ws.Range("N92").Value = -4465.8
$ws.Range("H122").Value = 1196.9667
$ws.Range("I122").Value = 612.25
$ws.Range("J122").Value = 1286.9231
$ws.Range("K122").Value = 5510.25
$ws.Range("L122").Value = 11582.3079
$ws.Range("M122").Value = -3060.25
$ws.Range("N122").Value = -16482.3079
$ws.Range("H129").Value = 1787.8889
$ws.Range("I129").Value = 933.25
$ws.Range("J129").Value = 2471.6
$ws.Range("K129").Value = 2799.75
$ws.Range("L129").Value = 7414.799999999999
$ws.Range("M129").Value = 2200.25
$ws.Range("N129").Value = -17414.8
$ws.Range("H131").Value = 1635224.5
$ws.Range("I131").Value = 2262922
$ws.Range("K131").Value = 6788766
$ws.Range("M131").Value = -6783726

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5727.421
$ws.Range("J70").Value = 5717.25
$ws.Range("L70").Value = 5717.25
$ws.Range("N70").Value = -6257.25
$ws.Range("H73").Value = 5727.421
$ws.Range("J73").Value = 5717.25
$ws.Range("L73").Value = 5717.25
$ws.Range("N73").Value = -7589.25
$ws.Range("H80").Value = 18464272
$ws.Range("J80").Value = 3152.7646
$ws.Range("L80").Value = 3152.7646
$ws.Range("N80").Value = -5148.7646
$ws.Range("H83").Value = 18464272
$ws.Range("J83").Value = 3152.7646
$ws.Range("L83").Value = 15763.823
$ws.Range("N83").Value = -25747.823
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4373.769
$ws.Range("J113").Value = 4682.3335
$ws.Range("K113").Value = 4373.769
$ws.Range("L113").Value = 4682.3335
$ws.Range("M113").Value = -2203.769
$ws.Range("N113").Value = -9022.333500000001
$ws.Range("H122").Value = 2454.2144
$ws.Range("I122").Value = 2342.889
$ws.Range("K122").Value = 7028.667
$ws.Range("M122").Value = -4578.667
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null
$ws.Range("H126").Value = 5246.5713
$ws.Range("I126").Value = 5062.263
$ws.Range("J126").Value = 6997.5
$ws.Range("K126").Value = 15186.789
$ws.Range("L126").Value = 20992.5
$ws.Range("M126").Value = -12716.789
$ws.Range("N126").Value = -25932.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6945179.5
$ws.Range("I16").Value = 7813202
$ws.Range("J16").Value = 998.5
$ws.Range("K16").Value = 7813202
$ws.Range("L16").Value = 998.5
$ws.Range("M16").Value = -7813032
$ws.Range("N16").Value = -1338.5
$ws.Range("H20").Value = 190
$ws.Range("I20").Value = 190
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 190
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 36
$ws.Range("N20").Value = $null
$ws.Range("H46").Value = 2405.1428
$ws.Range("I46").Value = 945.1429000000001
$ws.Range("K46").Value = 945.1429000000001
$ws.Range("M46").Value = -757.1429000000001
$ws.Range("H61").Value = 5932.857
$ws.Range("I61").Value = 6159.6
$ws.Range("K61").Value = 6159.6
$ws.Range("M61").Value = -5957.6
$ws.Range("H68").Value = 10420679
$ws.Range("I68").Value = 11496543
$ws.Range("K68").Value = 11496543
$ws.Range("M68").Value = -11495794
$ws.Range("H71").Value = 10420679
$ws.Range("I71").Value = 11496543
$ws.Range("K71").Value = 57482715
$ws.Range("M71").Value = -57478971
$ws.Range("H93").Value = 6898264.5
$ws.Range("I93").Value = 12501914
$ws.Range("K93").Value = 12501914
$ws.Range("M93").Value = -12500666
$ws.Range("H113").Value = 5932.857
$ws.Range("I113").Value = 6159.6
$ws.Range("K113").Value = 6159.6
$ws.Range("M113").Value = -3989.6
$ws.Range("H122").Value = 4692.8
$ws.Range("I122").Value = 4071.2917
$ws.Range("J122").Value = 19609
$ws.Range("K122").Value = 12213.8751
$ws.Range("L122").Value = 58827
$ws.Range("M122").Value = -9763.875100000001
$ws.Range("N122").Value = -63727
$ws.Range("H132").Value = 29503.695
$ws.Range("J132").Value = 2990.625
$ws.Range("L132").Value = 8971.875
$ws.Range("N132").Value = -14031.875
$ws.Range("H136").Value = 8186024.5
$ws.Range("I136").Value = 12003012
$ws.Range("J136").Value = 6764.4287
$ws.Range("K136").Value = 36009036
$ws.Range("L136").Value = 20293.2861
$ws.Range("M136").Value = -36006486
$ws.Range("N136").Value = -25393.2861

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 40124
$ws.Range("I40").Value = 43832
$ws.Range("J40").Value = 29000
$ws.Range("K40").Value = 43832
$ws.Range("L40").Value = 29000
$ws.Range("M40").Value = -43683
$ws.Range("N40").Value = -29298
$ws.Range("H81").Value = 7579944
$ws.Range("I81").Value = 11366841
$ws.Range("J81").Value = 6149.75
$ws.Range("K81").Value = 22733682
$ws.Range("L81").Value = 12299.5
$ws.Range("M81").Value = -22732621
$ws.Range("N81").Value = -14421.5
$ws.Range("H84").Value = 7579944
$ws.Range("I84").Value = 11366841
$ws.Range("J84").Value = 6149.75
$ws.Range("K84").Value = 113668410
$ws.Range("L84").Value = 61497.5
$ws.Range("M84").Value = -113663106
$ws.Range("N84").Value = -72105.5
$ws.Range("H107").Value = 801.12
$ws.Range("I107").Value = 785
$ws.Range("J107").Value = 852.1667
$ws.Range("K107").Value = 2355
$ws.Range("L107").Value = 2556.5001
$ws.Range("M107").Value = -435
$ws.Range("N107").Value = -6396.5001
$ws.Range("H113").Value = 2434.4666
$ws.Range("I113").Value = 979.1429000000001
$ws.Range("J113").Value = 3707.875
$ws.Range("K113").Value = 2937.4287
$ws.Range("L113").Value = 11123.625
$ws.Range("M113").Value = -767.4287000000004
$ws.Range("N113").Value = -15463.625
$ws.Range("H128").Value = 149750
$ws.Range("J128").Value = 149750
$ws.Range("L128").Value = 149750
$ws.Range("N128").Value = -159710
$ws.Range("H132").Value = 2713.9
$ws.Range("I132").Value = 2312.0833
$ws.Range("K132").Value = 6936.249899999999
$ws.Range("M132").Value = -4406.249899999999
$ws.Range("H136").Value = 1672.9436
$ws.Range("I136").Value = 1838.8704
$ws.Range("K136").Value = 5516.6112
$ws.Range("M136").Value = -2966.6112
